$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''38.269.71'
$ws.Range("E2").Value = '  +3.36%  '

# Row 3
$ws.Range("D3").Value = '''2.064.67'
$ws.Range("E3").Value = '  +2.77%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").Value = '''230.62'
$ws.Range("E5").Value = '  +2.14%  '

# Row 6
$ws.Range("D6").Value = '''0.617'
$ws.Range("E6").Value = '  +2.08%  '

# Row 7
$ws.Range("D7").Value = '''61.05'
$ws.Range("E7").Value = '  +10.67%  '

# Row 8
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("D9").Value = '''0.388'
$ws.Range("E9").Value = '  +4.23%  '

# Row 10
$ws.Range("D10").Value = '''0.0813'
$ws.Range("E10").Value = '  +4.72%  '

# Row 11
$ws.Range("E11").Value = '  +2.07%  '

# Row 12
$ws.Range("D12").Value = '''14.85'
$ws.Range("E12").Value = '  +5.86%  '

# Row 13
$ws.Range("D13").Value = '''2.367.22'
$ws.Range("E13").Value = '  +2.54%  '

# Row 14
$ws.Range("D14").Value = '''21.45'
$ws.Range("E14").Value = '  +8.76%  '

# Row 15
$ws.Range("D15").Value = '''0.765'
$ws.Range("E15").Value = '  +4.28%  '

# Row 16
$ws.Range("D16").Value = '''5.34'
$ws.Range("E16").Value = '  +3.79%  '

# Row 17
$ws.Range("D17").Value = '''2.067.18'
$ws.Range("E17").Value = '  +2.87%  '

# Row 18
$ws.Range("D18").Value = '''38.168.88'
$ws.Range("E18").Value = '  +3.29%  '

# Row 19
$ws.Range("E19").Value = '  +2.59%  '

# Row 20
$ws.Range("D20").Value = '''70.19'
$ws.Range("E20").Value = '  +2.90%  '

# Row 21
$ws.Range("D21").Value = '''0.0₃0837'
$ws.Range("E21").Value = '  +3.31%  '

# Row 22
$ws.Range("D22").Value = '''226.06'
$ws.Range("E22").Value = '  +1.49%  '

# Row 23
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.04%  '

# Row 24
$ws.Range("D24").Value = '''2.43'
$ws.Range("E24").Value = '  -0.03%  '

# Row 25
$ws.Range("E25").Value = '  +3.67%  '

# Row 26
$ws.Range("D26").Value = '''9.33'
$ws.Range("E26").Value = '  +4.55%  '

# Row 27
$ws.Range("D27").Value = '''166.23'
$ws.Range("E27").Value = '  +1.03%  '

# Row 28
$ws.Range("D28").Value = '''0.134'
$ws.Range("E28").Value = '  +6.52%  '

# Row 29
$ws.Range("D29").Value = '''19.15'
$ws.Range("E29").Value = '  +3.28%  '

# Row 30
$ws.Range("E30").Value = '  +2.10%  '

# Row 31
$ws.Range("E31").Value = '  +2.93%  '

# Row 32
$ws.Range("D32").Value = '''4.57'
$ws.Range("E32").Value = '  +4.37%  '

# Row 33
$ws.Range("D33").Value = '''4.64'
$ws.Range("E33").Value = '  +4.62%  '

# Row 34
$ws.Range("E34").Value = '  +9.99%  '

# Row 35
$ws.Range("D35").Value = '''0.0608'
$ws.Range("E35").Value = '  +1.47%  '

# Row 36
$ws.Range("D36").Value = '''2.32'
$ws.Range("E36").Value = '  +0.22%  '

# Row 37
$ws.Range("D37").Value = '''6.23'
$ws.Range("E37").Value = '  +16.66%  '

# Row 38
$ws.Range("D38").Value = '''3.33'
$ws.Range("E38").Value = '  +6.25%  '

# Row 39
$ws.Range("E39").Value = '  -0.01%  '

# Row 40
$ws.Range("D40").Value = '''1.529.68'
$ws.Range("E40").Value = '  +4.58%  '

# Row 41
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '''17.18'
$ws.Range("E41").Value = '  +8.43%  '

# Row 42
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '''98.35'
$ws.Range("E42").Value = '  +4.30%  '

# Row 43
$ws.Range("D43").Value = '''0.0218'
$ws.Range("E43").Value = '  +3.14%  '

# Row 44
$ws.Range("E44").Value = '  +4.14%  '

# Row 45
$ws.Range("D45").Value = '''0.0930'
$ws.Range("E45").Value = '  +2.12%  '

# Row 46
$ws.Range("E46").Value = '  +1.55%  '

# Row 47
$ws.Range("E47").Value = '  -4.40%  '

# Row 48
$ws.Range("D48").Value = '''1.03'
$ws.Range("E48").Value = '  +3.17%  '

# Row 49
$ws.Range("D49").Value = '''2.96'
$ws.Range("E49").Value = '  +2.64%  '

# Row 50
$ws.Range("D50").Value = '''7.12'
$ws.Range("E50").Value = '  +1.14%  '

# Row 51
$ws.Range("D51").Value = '''2.256.19'
$ws.Range("E51").Value = '  +2.60%  '
